$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the betting-data columns (F:V) between rows 49/51 and 50/52 ---
# (indices / dates in A:E stay put; only the match details moved rows.)
$r49 = $ws.Range("F49:V49").Value2
$r51 = $ws.Range("F51:V51").Value2
$ws.Range("F49:V49").Value2 = $r51
$ws.Range("F51:V51").Value2 = $r49

$r50 = $ws.Range("F50:V50").Value2
$r52 = $ws.Range("F52:V52").Value2
$ws.Range("F50:V50").Value2 = $r52
$ws.Range("F52:V52").Value2 = $r50

# --- 2) Swap the betting-data columns (F:V) between rows 91 and 93 ---
$r91 = $ws.Range("F91:V91").Value2
$r93 = $ws.Range("F93:V93").Value2
$ws.Range("F91:V91").Value2 = $r93
$ws.Range("F93:V93").Value2 = $r91

# --- 3) Append four new match rows (136-139), matching the formatting of
#        the previous last row (135: bold/bordered index, dated kickoff) ---
$ws.Range("A135:V135").Copy()
$ws.Range("A136:V139").PasteSpecial(-4122)

$newRows = @(
    @(135, "poland", "division-2", "2023-2024", 45234.5625, "KKS Kalisz", 2, "Ol. Grudziadz", 2, 1.75, "03/11/2023 01:43", 2.35, "04/11/2023 13:16", 3.52, "03/11/2023 01:43", 3.44, "04/11/2023 13:16", 3.97, "03/11/2023 01:43", 2.82, "04/11/2023 13:16", "https://www.betexplorer.com/football/poland/division-2/kks-kalisz-ol-grudziadz/nsMDXJ54/"),
    @(136, "poland", "division-2", "2023-2024", 45234.5625, "Sandecja Nowy S.", 0, "Lech Poznan II", 0, 1.78, "03/11/2023 01:43", 1.85, "04/11/2023 13:23", 3.59, "03/11/2023 01:43", 3.62, "04/11/2023 13:23", 3.66, "03/11/2023 01:43", 3.93, "04/11/2023 13:23", "https://www.betexplorer.com/football/poland/division-2/sandecja-nowy-s-lech-poznan/KfKLVczH/"),
    @(137, "poland", "division-2", "2023-2024", 45234.69791666666, "Chojniczanka", 1, "Pogon Siedlce", 2, 1.99, "03/11/2023 05:12", 2.2, "04/11/2023 16:34", 3.26, "03/11/2023 05:12", 3.27, "04/11/2023 16:34", 3.26, "03/11/2023 05:12", 3.23, "04/11/2023 16:34", "https://www.betexplorer.com/football/poland/division-2/chojniczanka-pogon-siedlce/OSliGIzU/"),
    @(138, "poland", "division-2", "2023-2024", 45234.80208333334, "Stomil Olsztyn", 2, "Skra", 1, 2.31, "03/11/2023 07:43", 2.58, "04/11/2023 19:14", 3.02, "03/11/2023 07:43", 3.01, "04/11/2023 19:14", 2.85, "03/11/2023 07:43", 2.84, "04/11/2023 19:14", "https://www.betexplorer.com/football/poland/division-2/stomil-olsztyn-skra-czestochowa/z7I9Yajb/")
)

$startRow = 136
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $newRows[$i]
    $arr = New-Object 'object[,]' 1,22
    for ($c = 0; $c -lt 22; $c++) {
        $arr[0,$c] = $row[$c]
    }
    $targetRow = $startRow + $i
    $ws.Range("A$($targetRow):V$($targetRow)").Value2 = $arr
}
